# Update "想去人数" (want-to-go count) figures in column F across the
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types) sheets,
# matching the freshly generated data snapshot.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll  = $wb.Worksheets.Item("全部类型")

# column F is the 6th column
$col = 6

# --- 展览 sheet ---
$expoUpdates = @{
    2  = 640
    4  = 50
    6  = 418
    7  = 241
    8  = 13271
    9  = 55
    10 = 46
    11 = 5374
    12 = 560
    13 = 28
    15 = 42
    17 = 52
    19 = 708
    20 = 2877
    21 = 7159
    22 = 1171
    23 = 3658
    24 = 226
    25 = 55
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Cells.Item($row, $col).Value = $expoUpdates[$row]
}

# --- 演出 sheet ---
$showUpdates = @{
    2 = 32
}
foreach ($row in $showUpdates.Keys) {
    $wsShow.Cells.Item($row, $col).Value = $showUpdates[$row]
}

# --- 全部类型 sheet ---
$allUpdates = @{
    2  = 640
    4  = 50
    6  = 32
    7  = 418
    8  = 241
    9  = 13271
    10 = 55
    11 = 46
    12 = 5374
    13 = 560
    14 = 28
    16 = 42
    18 = 52
    20 = 708
    21 = 2877
    23 = 7159
    24 = 1171
    25 = 3658
    26 = 226
    27 = 55
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, $col).Value = $allUpdates[$row]
}
